$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Profiles")

# Update column A timestamps for rows 4-19 so they land on exact
# hourly boundaries relative to the day start in A3 (44410).
for ($i = 4; $i -le 19; $i++) {
    $ws.Cells.Item($i, 1).Value = 44410 + (($i - 3) / 24)
}

# Remove the now-unused rows 20-75 (profiles shortened to 16h/19 rows).
$ws.Range("A20:C75").EntireRow.Delete()
